# Apply the "What we like" / "What we don't like" / meta-description edits
# described by the commit "Added many more features".

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Visually stunning graphics and game symbols" "Visually stunning graphics"
Replace-Text "Cascading reels and free spins bonus increase chances of winning" "Unique cascading reels and free spins bonus"
Replace-Text "Wild symbol substitutes for all regular symbols" "Wild symbol for easy winning combinations"
Replace-Text "Flexible betting range suitable for all players" "Wide range of betting options"
Replace-Text "Limited bonus features compared to other slot games" "None"
Replace-Text "Audio effects can be repetitive" "None"
Replace-Text "Read our review of Invaders from the Planet Moolah and play for free. Enjoy cascading reels and flexible betting range with this space-themed slot game." "Read our review of Invaders from the Planet Moolah and play this game for free to experience its unique features."
